# Insert two new weekly price rows (Primera / Segunda) for "Pepino ensalada"
# at the top of the date-ordered block starting at row 252. All existing
# rows from 252 downward shift down by two rows (to 254.. and the former
# 353/354 become the new last rows 355/356).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 252 and below down by two rows, carrying formatting along
# (mirrors the row above, same as a normal Excel row insert).
$ws.Rows("252:253").Insert()

# --- New row 252: Calidad "Primera" ---
$ws.Range("A252").Value = 1
$ws.Range("B252").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C252").Value = "Arica y Parinacota"
$ws.Range("D252").Value = 44795
$ws.Range("E252").Value = 15
$ws.Range("F252").Value = 100112043
$ws.Range("G252").Value = "Pepino ensalada"
$ws.Range("H252").Value = "Sin especificar"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 150
$ws.Range("K252").Value = 18000
$ws.Range("L252").Value = 19000
$ws.Range("M252").Value = 18500
$ws.Range("N252").Value = "$/caja 70 unidades"
$ws.Range("O252").Value = "Región de Arica y Parinacota"
$ws.Range("P252").Value = 264
$ws.Range("Q252").Value = 70
$ws.Range("R252").Value = "Hortaliza"

# --- New row 253: Calidad "Segunda" ---
$ws.Range("A253").Value = 1
$ws.Range("B253").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C253").Value = "Arica y Parinacota"
$ws.Range("D253").Value = 44795
$ws.Range("E253").Value = 15
$ws.Range("F253").Value = 100112043
$ws.Range("G253").Value = "Pepino ensalada"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Segunda"
$ws.Range("J253").Value = 200
$ws.Range("K253").Value = 14000
$ws.Range("L253").Value = 15000
$ws.Range("M253").Value = 14500
$ws.Range("N253").Value = "$/caja 100 unidades"
$ws.Range("O253").Value = "Región de Arica y Parinacota"
$ws.Range("P253").Value = 145
$ws.Range("Q253").Value = 100
$ws.Range("R253").Value = "Hortaliza"
